$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns before column D (shifts existing D:K data to F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formatting from column F (the old column D, now shifted) into new D:E columns
$ws.Range("F1:F102").Copy()
$ws.Range("D1:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New quarterly data for the newly inserted columns D (most recent) and E
$data = @{
    7 = @(43465, 43373)
    8 = @(131700, 123900)
    9 = @(81800, 79200)
    10 = @(49900, 44700)
    11 = @($null, $null)
    12 = @("NA", "NA")
    13 = @(0, 0)
    14 = @(500, 500)
    15 = @(5300, 5000)
    16 = @($null, $null)
    17 = @(120300, 114700)
    18 = @(11400, 9200)
    19 = @($null, $null)
    20 = @(-1800, -800)
    21 = @(14900, 13400)
    22 = @(0, 0)
    23 = @(9600, 8400)
    24 = @(2100, 2100)
    25 = @(0, 0)
    26 = @(7500, 6300)
    27 = @(7500, 6300)
    28 = @(0, 0)
    29 = @("NA", "NA")
    30 = @(0, 0)
    31 = @(0, 0)
    32 = @(1800, 800)
    33 = @(7500, 6300)
    34 = @(0, 0)
    35 = @(7500, 6300)
    38 = @(43465, 43373)
    39 = @($null, $null)
    40 = @($null, $null)
    41 = @(45000, 44900)
    42 = @(0, 0)
    43 = @(122400, 109800)
    44 = @(0, 0)
    45 = @(10400, 6500)
    46 = @(177800, 161200)
    47 = @(0, 0)
    48 = @(6700, 6600)
    49 = @(376100, 371800)
    50 = @(0, 0)
    51 = @(0, 0)
    52 = @(10000, 9700)
    53 = @(0, 0)
    54 = @(570500, 549300)
    55 = @($null, $null)
    56 = @($null, $null)
    57 = @(24400, 10800)
    58 = @(0, 0)
    59 = @(50400, 44200)
    60 = @(74800, 55000)
    61 = @(120100, 119000)
    62 = @(22000, 20300)
    63 = @(0, 0)
    64 = @(0, 0)
    65 = @(0, 0)
    66 = @(216900, 194300)
    67 = @($null, $null)
    68 = @(0, 0)
    69 = @(0, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(152700, 145200)
    73 = @(0, 0)
    74 = @(0, 0)
    75 = @(0, 0)
    76 = @(353700, 355000)
    77 = @(0, 0)
    80 = @(43465, 43373)
    81 = @(7500, 6300)
    82 = @($null, $null)
    83 = @(5300, 5000)
    84 = @(0, 0)
    85 = @(0, 0)
    86 = @(0, 0)
    87 = @(0, 0)
    88 = @(0, 0)
    89 = @(25800, 13200)
    90 = @($null, $null)
    91 = @(-900, -1600)
    92 = @(0, 0)
    93 = @(0, 0)
    94 = @(-6400, -11700)
    95 = @($null, $null)
    96 = @(0, 0)
    97 = @(0, 0)
    98 = @(0, 0)
    99 = @(0, 0)
    100 = @(-19300, 33300)
    101 = @(0, -200)
    102 = @(0, 34600)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    if ($null -ne $vals[0]) { $ws.Cells.Item($row, 4).Value = $vals[0] }
    if ($null -ne $vals[1]) { $ws.Cells.Item($row, 5).Value = $vals[1] }
}